# Generate Report for handoff
#
# The "09caa6ba-7754-4c23-b998-8e4b62501455" row on both the zh-cn and
# de-de localization-status sheets was re-handed-off; its
# "Latest Handoff Datetime" (column D) needs to move forward while the
# neighboring "656b4403-eba8-4361-af3f-dc3c9dfc61f7" row's handoff
# timestamp is untouched.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 4 is the 09caa6ba entry, row 5 is the 656b4403 entry.
$wsZh.Range("D4").Value = "2016-01-18 02:28:32"

# de-de sheet: row 4 is the 09caa6ba entry, row 5 is the 656b4403 entry.
$wsDe.Range("D4").Value = "2016-01-18 02:28:44"
